# Apply updated Betfair Back/Lay odds values for 2026-01-07 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 5.2  # I2
$ws.Cells.Item(2, 12).Value = 1.5  # L2
$ws.Cells.Item(2, 20).Value = 2.22  # T2
$ws.Cells.Item(2, 22).Value = 1.23  # V2
$ws.Cells.Item(3, 16).Value = 2.28  # P3
$ws.Cells.Item(4, 6).Value = 4.8  # F4
$ws.Cells.Item(4, 7).Value = 5.8  # G4
$ws.Cells.Item(4, 9).Value = 1.96  # I4
$ws.Cells.Item(4, 16).Value = 1.7  # P4
$ws.Cells.Item(4, 17).Value = 2.18  # Q4
$ws.Cells.Item(4, 24).Value = 13.5  # X4
$ws.Cells.Item(4, 25).Value = 8.6  # Y4
$ws.Cells.Item(4, 26).Value = 12.5  # Z4
$ws.Cells.Item(4, 27).Value = 22  # AA4
$ws.Cells.Item(4, 28).Value = 18  # AB4
$ws.Cells.Item(4, 29).Value = 9.6  # AC4
$ws.Cells.Item(4, 31).Value = 23  # AE4
$ws.Cells.Item(4, 32).Value = 42  # AF4
$ws.Cells.Item(4, 33).Value = 25  # AG4
$ws.Cells.Item(4, 34).Value = 26  # AH4
$ws.Cells.Item(4, 41).Value = 17  # AO4
$ws.Cells.Item(5, 6).Value = 3.45  # F5
$ws.Cells.Item(5, 7).Value = 36  # G5
$ws.Cells.Item(5, 8).Value = 1.09  # H5
$ws.Cells.Item(5, 9).Value = 2.14  # I5
$ws.Cells.Item(5, 10).Value = 2  # J5
$ws.Cells.Item(5, 11).Value = 32  # K5
$ws.Cells.Item(5, 14).Value = 1.1  # N5
$ws.Cells.Item(5, 15).Value = 1.15  # O5
$ws.Cells.Item(5, 16).Value = 1.24  # P5
$ws.Cells.Item(5, 17).Value = 1.15  # Q5
$ws.Cells.Item(5, 19).Value = 1.28  # S5
$ws.Cells.Item(5, 22).Value = 1.88  # V5
$ws.Cells.Item(5, 23).Value = 1.03  # W5
$ws.Cells.Item(7, 7).Value = 17  # G7
$ws.Cells.Item(7, 23).Value = 1.06  # W7
$ws.Cells.Item(8, 10).Value = 5.4  # J8
$ws.Cells.Item(8, 11).Value = 5.5  # K8
$ws.Cells.Item(8, 24).Value = 15  # X8
$ws.Cells.Item(8, 25).Value = 28  # Y8
$ws.Cells.Item(8, 34).Value = 40  # AH8
$ws.Cells.Item(9, 7).Value = 2.98  # G9
$ws.Cells.Item(9, 23).Value = 1.5  # W9
$ws.Cells.Item(9, 39).Value = 95  # AM9
$ws.Cells.Item(9, 40).Value = 32  # AN9
$ws.Cells.Item(9, 41).Value = 27  # AO9
$ws.Cells.Item(10, 8).Value = 2.16  # H10
$ws.Cells.Item(10, 17).Value = 1.85  # Q10
$ws.Cells.Item(10, 22).Value = 1.84  # V10
$ws.Cells.Item(10, 24).Value = 17  # X10
$ws.Cells.Item(10, 40).Value = 34  # AN10
$ws.Cells.Item(11, 7).Value = 1.92  # G11
$ws.Cells.Item(11, 8).Value = 4.7  # H11
$ws.Cells.Item(11, 16).Value = 1.88  # P11
$ws.Cells.Item(11, 23).Value = 2.08  # W11
$ws.Cells.Item(11, 28).Value = 8.4  # AB11
$ws.Cells.Item(12, 8).Value = 2.36  # H12
$ws.Cells.Item(12, 9).Value = 2.38  # I12
$ws.Cells.Item(12, 16).Value = 2  # P12
$ws.Cells.Item(12, 17).Value = 1.97  # Q12
$ws.Cells.Item(12, 18).Value = 1.38  # R12
$ws.Cells.Item(12, 20).Value = 1.76  # T12
$ws.Cells.Item(12, 22).Value = 1.72  # V12
$ws.Cells.Item(12, 31).Value = 23  # AE12
$ws.Cells.Item(12, 33).Value = 14  # AG12
$ws.Cells.Item(13, 12).Value = 1.36  # L13
$ws.Cells.Item(13, 15).Value = 1.27  # O13
$ws.Cells.Item(13, 16).Value = 2.14  # P13
$ws.Cells.Item(13, 18).Value = 1.45  # R13
$ws.Cells.Item(13, 19).Value = 3.1  # S13
$ws.Cells.Item(13, 20).Value = 1.7  # T13
$ws.Cells.Item(13, 21).Value = 2.36  # U13
$ws.Cells.Item(13, 25).Value = 15  # Y13
$ws.Cells.Item(13, 28).Value = 11.5  # AB13
$ws.Cells.Item(13, 29).Value = 8  # AC13
$ws.Cells.Item(13, 35).Value = 44  # AI13
$ws.Cells.Item(13, 37).Value = 22  # AK13
$ws.Cells.Item(14, 6).Value = 1.44  # F14
$ws.Cells.Item(14, 7).Value = 1.45  # G14
$ws.Cells.Item(14, 8).Value = 7.6  # H14
$ws.Cells.Item(14, 9).Value = 7.8  # I14
$ws.Cells.Item(14, 10).Value = 5.5  # J14
$ws.Cells.Item(14, 18).Value = 1.84  # R14
$ws.Cells.Item(14, 20).Value = 1.65  # T14
$ws.Cells.Item(14, 22).Value = 1.14  # V14
$ws.Cells.Item(14, 23).Value = 3.2  # W14
$ws.Cells.Item(14, 26).Value = 75  # Z14
$ws.Cells.Item(14, 30).Value = 28  # AD14
$ws.Cells.Item(14, 31).Value = 85  # AE14
$ws.Cells.Item(14, 32).Value = 12  # AF14
$ws.Cells.Item(14, 33).Value = 10  # AG14
$ws.Cells.Item(14, 36).Value = 14  # AJ14
$ws.Cells.Item(15, 22).Value = 1.23  # V15
$ws.Cells.Item(16, 6).Value = 2.56  # F16
$ws.Cells.Item(16, 7).Value = 2.58  # G16
$ws.Cells.Item(16, 16).Value = 1.65  # P16
$ws.Cells.Item(16, 17).Value = 2.46  # Q16
$ws.Cells.Item(16, 21).Value = 1.91  # U16
$ws.Cells.Item(16, 23).Value = 1.63  # W16
$ws.Cells.Item(16, 29).Value = 7  # AC16
$ws.Cells.Item(16, 32).Value = 14  # AF16
$ws.Cells.Item(16, 36).Value = 36  # AJ16
$ws.Cells.Item(16, 37).Value = 32  # AK16
$ws.Cells.Item(16, 40).Value = 32  # AN16
$ws.Cells.Item(17, 14).Value = 3  # N17
$ws.Cells.Item(17, 24).Value = 9.4  # X17
$ws.Cells.Item(18, 6).Value = 9.6  # F18
$ws.Cells.Item(18, 14).Value = 5  # N18
$ws.Cells.Item(18, 16).Value = 2.38  # P18
$ws.Cells.Item(18, 17).Value = 1.69  # Q18
$ws.Cells.Item(18, 18).Value = 1.54  # R18
$ws.Cells.Item(18, 26).Value = 8.4  # Z18
$ws.Cells.Item(19, 6).Value = 5.3  # F19
$ws.Cells.Item(19, 8).Value = 1.74  # H19
$ws.Cells.Item(19, 9).Value = 1.75  # I19
$ws.Cells.Item(19, 10).Value = 4.1  # J19
$ws.Cells.Item(19, 16).Value = 2.18  # P19
$ws.Cells.Item(19, 19).Value = 3  # S19
$ws.Cells.Item(19, 22).Value = 2.32  # V19
$ws.Cells.Item(19, 27).Value = 17.5  # AA19
$ws.Cells.Item(19, 28).Value = 20  # AB19
$ws.Cells.Item(19, 33).Value = 19.5  # AG19
$ws.Cells.Item(19, 35).Value = 30  # AI19
$ws.Cells.Item(19, 38).Value = 70  # AL19
$ws.Cells.Item(19, 41).Value = 9  # AO19
$ws.Cells.Item(20, 6).Value = 1.77  # F20
$ws.Cells.Item(20, 7).Value = 1.78  # G20
$ws.Cells.Item(20, 16).Value = 2.12  # P20
$ws.Cells.Item(20, 17).Value = 1.87  # Q20
$ws.Cells.Item(20, 18).Value = 1.43  # R20
$ws.Cells.Item(20, 19).Value = 3.2  # S20
$ws.Cells.Item(20, 20).Value = 1.85  # T20
$ws.Cells.Item(20, 29).Value = 8.4  # AC20
$ws.Cells.Item(20, 34).Value = 18  # AH20
$ws.Cells.Item(20, 36).Value = 17  # AJ20
$ws.Cells.Item(20, 41).Value = 65  # AO20
$ws.Cells.Item(21, 10).Value = 1.04  # J21
